$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.562.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.97%  "

$ws.Range("D3").Value = "'2.223.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.18%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'229.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "

$ws.Range("D6").Value = "'0.623"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("D7").Value = "'61.26"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.70%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  +3.36%  "

$ws.Range("D10").Value = "'58.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("D11").Value = "'0.0882"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.64%  "

$ws.Range("E12").Value = "  +0.22%  "

$ws.Range("D13").Value = "'2.556.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.33%  "

$ws.Range("D14").Value = "'15.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.68%  "

$ws.Range("D15").Value = "'21.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.89%  "

$ws.Range("D16").Value = "'0.795"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.86%  "

$ws.Range("D17").Value = "'5.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.33%  "

$ws.Range("D18").Value = "'2.222.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.02%  "

$ws.Range("D19").Value = "'41.534.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.02%  "

$ws.Range("D20").Value = "'72.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.62%  "

$ws.Range("D21").Value = "'0.0₃0891"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.74%  "

$ws.Range("E22").Value = "  -0.34%  "

$ws.Range("D23").Value = "'249.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.90%  "

$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").Value = "'2.37"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").Value = "'2.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.33%  "

$ws.Range("D27").Value = "'9.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.52%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.142"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.53%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'167.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.17%  "

$ws.Range("D30").Value = "'19.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.50%  "

$ws.Range("D31").Value = "'1.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.71%  "

$ws.Range("D32").Value = "'2.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.32%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").Value = "'4.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.44%  "

$ws.Range("E35").Value = "  +1.23%  "

$ws.Range("D36").Value = "'0.0622"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.25%  "

$ws.Range("D37").Value = "'6.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.92%  "

$ws.Range("E38").Value = "  +0.74%  "

$ws.Range("D39").Value = "'2.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.10%  "

$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.06%  "

$ws.Range("E41").Value = "  +28.16%  "

$ws.Range("D42").Value = "'4.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.64%  "

$ws.Range("D44").Value = "'8.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.33%  "

$ws.Range("D45").Value = "'0.0983"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.40%  "

$ws.Range("D46").Value = "'99.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.37%  "

$ws.Range("D47").Value = "'1.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.92%  "

$ws.Range("D48").Value = "'1.465.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.30%  "

$ws.Range("D49").Value = "'16.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.18%  "

$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("E51").Value = "  -0.98%  "
